$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I ("Other found locations")
$ws.Range("I1").Value = 'Other found locations'

# Row 2
$ws.Range("F2").Value = 'not found'
$ws.Range("G2").Value = 'N/A'
$ws.Range("I2").Value = ''

# Row 3
$ws.Range("I3").Value = ''

# Row 4
$ws.Range("E4").Value = '[Barret%BP%coreGivesNoEmail%1,  Bartels%CL%coreGivesNoEmail%1,  Bauer%R%coreGivesNoEmail%1,  Brevoort%P%coreGivesNoEmail%1,  Brinkeborn%RM%coreGivesNoEmail%1,  Carr%RJ%coreGivesNoEmail%1,  Ernst%E%coreGivesNoEmail%1,  Giles%JT%coreGivesNoEmail%1,  Hoheisel%D%coreGivesNoEmail%1,  Management%of Influenza in the Southern Hemisphere Trialists Study Group%coreGivesNoEmail%1,  Muller-Jakic%B%coreGivesNoEmail%1,  Stimpel%M%coreGivesNoEmail%1,  Turner%RB%coreGivesNoEmail%1,  Winther%B%coreGivesNoEmail%1]'
$ws.Range("F4").Value = 'not found'
$ws.Range("G4").Value = 'N/A'
$ws.Range("I4").Value = ''

# Row 5
$ws.Range("F5").Value = 'not found'
$ws.Range("G5").Value = 'N/A'
$ws.Range("I5").Value = ''

# Row 6
$ws.Range("I6").Value = ''

# Row 7
$ws.Range("I7").Value = ''

# Row 8
$ws.Range("E8").Value = '[ Wolfram%Grimm%null%2,    Hans-Helge%Müller%null%1,  Wolfram%Grimm%null%0,  Hans-Helge%Müller%null%1]'
$ws.Range("I8").Value = ''

# Row 9
$ws.Range("I9").Value = ''

# Row 10
$ws.Range("I10").Value = ''

# Row 11
$ws.Range("E11").Value = '[Joelle%O’Neil%NULL%1,   Susan%Hughes%susan.hughes@fresno.ucsf.edu%1,   Andrea%Lourie%NULL%1,   John%Zweifler%NULL%1]'
$ws.Range("I11").Value = '_PMC_elsevier'

# Row 12
$ws.Range("E12").Value = '[M.%Jawad%NULL%1,   R.%Schoop%NULL%1,   A.%Suter%NULL%1,   P.%Klein%NULL%1,   R.%Eccles%NULL%1]'
$ws.Range("I12").Value = '_PMC'

# Row 13
$ws.Range("E13").Value = '[E.%Tiralongo%NULL%1,   R. A.%Lea%NULL%1,   S. S.%Wee%NULL%1,   M. M.%Hanna%NULL%1,   L. R.%Griffiths%NULL%1]'
$ws.Range("I13").Value = '_PMC'

# Row 14
$ws.Range("I14").Value = ''

# Row 15
$ws.Range("E15").Value = '[Steven J.%Sperber%ssperber@humed.com%1,   Leena P.%Shah%NULL%1,   Richard D.%Gilbert%NULL%1,   Thomas W.%Ritchey%NULL%1,   Arnold S.%Monto%NULL%1]'
$ws.Range("I15").Value = '_PMC'

# Row 16
$ws.Range("F16").Value = 'not found'
$ws.Range("G16").Value = 'N/A'
$ws.Range("I16").Value = ''

# Row 17
$ws.Range("I17").Value = ''

# Row 18
$ws.Range("E18").Value = '[ M.%Dorn%null%2,    E.%Knick%null%1,    G.%Lewith%null%1,  M.%Dorn%null%0,  E.%Knick%null%1,  G.%Lewith%null%1]'
$ws.Range("I18").Value = ''
